$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.484.39'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '3.815.32'
$ws.Range("E3").Value = '  +3.26%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '410.26'
$ws.Range("E5").Value = '  -2.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.80'
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("D7").Value = '3.804.99'
$ws.Range("E7").Value = '  +3.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.610'
$ws.Range("E8").Value = '  -5.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.732'
$ws.Range("E10").Value = '  -5.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.169'
$ws.Range("E11").Value = '  -7.25%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000363'
$ws.Range("E12").Value = '  -8.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.86'
$ws.Range("E13").Value = '  -5.42%  '
$ws.Range("D14").Value = '4.407.89'
$ws.Range("E14").Value = '  +2.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '10.03'
$ws.Range("E15").Value = '  -5.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.68'
$ws.Range("E16").Value = '  +17.80%  '
$ws.Range("E17").Value = '  -1.33%  '
$ws.Range("D18").Value = '3.822.69'
$ws.Range("E18").Value = '  +3.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '19.55'
$ws.Range("E19").Value = '  -5.70%  '
$ws.Range("D20").Value = '66.881.82'
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("E21").Value = '  -6.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '410.60'
$ws.Range("E22").Value = '  -7.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.36'
$ws.Range("E23").Value = '  -12.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.10'
$ws.Range("E24").Value = '  -5.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.05'
$ws.Range("E25").Value = '  -3.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '36.40'
$ws.Range("E26").Value = '  -2.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.71'
$ws.Range("E27").Value = '  +13.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.14'
$ws.Range("E28").Value = '  -6.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.44'
$ws.Range("E29").Value = '  -7.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '689.95'
$ws.Range("E30").Value = '  +5.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.45'
$ws.Range("E31").Value = '  -2.13%  '
$ws.Range("E32").Value = '  -2.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.73'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.15'
$ws.Range("E34").Value = '  -2.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.152'
$ws.Range("E35").Value = '  -7.82%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.64'
$ws.Range("E36").Value = '  -7.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = '0.0₃0792'
$ws.Range("E38").Value = '  +7.62%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '54.72'
$ws.Range("E39").Value = '  -4.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.18'
$ws.Range("E40").Value = '  +1.86%  '
$ws.Range("E41").Value = '  -7.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.996'
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '149.76'
$ws.Range("E43").Value = '  +0.44%  '
$ws.Range("E44").Value = '  -9.48%  '
$ws.Range("E45").Value = '  +2.06%  '
$ws.Range("B46").Value = 'LidoDAOToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.30'
$ws.Range("E46").Value = '  -4.20%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.51'
$ws.Range("E47").Value = '  -10.32%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.07'
$ws.Range("E48").Value = '  -2.09%  '
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.10'
$ws.Range("E49").Value = '  -4.97%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.78'
$ws.Range("E50").Value = '  -4.06%  '
$ws.Range("B51").Value = 'WEMIXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.55'
$ws.Range("E51").Value = '  -3.87%  '
